# This script applies a cyclic rotation of the per-observation data among
# rows 2, 3, 5, 6 and 7 of the active worksheet, as described by the
# commit diff (the location/meta columns already hold equal values across
# these rows, so only the columns below actually show up as changed):
#
#   new row 2 <= old row 7
#   new row 3 <= old row 5
#   new row 5 <= old row 2
#   new row 6 <= old row 3
#   new row 7 <= old row 6
#
# Columns that rotate: A (Id), B (Taxonsorteringsordning), E (TaxonId),
# F (Artnamn), G (Vetenskapligt namn), H (Auktor), I (Antal), J (Enhet),
# L (empty "Kon" placeholder cell), Q (Ost) and R (Nord).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextOrBlank {
    param(
        [object]$cell,
        [object]$value,
        [bool]$forceText
    )
    # Whatever was there before, first drop the contents...
    $cell.ClearContents()
    if ($value -eq $null) {
        # ...then re-materialise an empty placeholder cell (mirrors the
        # blank inlineStr cells already present in the workbook).
        $cell.NumberFormat = "General"
    } else {
        if ($forceText) {
            # Force text storage so a numeric-looking value (e.g. "5") is
            # kept as text instead of being coerced into a number.
            $cell.NumberFormat = "@"
        }
        $cell.Value2 = $value
    }
}

function Set-RowData {
    param(
        [int]$row,
        [double]$a,
        [double]$b,
        [double]$e,
        [string]$f,
        [string]$g,
        [string]$h,
        [object]$antal,
        [object]$enhet,
        [bool]$hasL,
        [double]$q,
        [double]$r
    )

    $ws.Cells.Item($row, 1).Value2 = $a    # A - Id
    $ws.Cells.Item($row, 2).Value2 = $b    # B - Taxonsorteringsordning
    $ws.Cells.Item($row, 5).Value2 = $e    # E - TaxonId
    $ws.Cells.Item($row, 6).Value2 = $f    # F - Artnamn
    $ws.Cells.Item($row, 7).Value2 = $g    # G - Vetenskapligt namn
    $ws.Cells.Item($row, 8).Value2 = $h    # H - Auktor

    Set-TextOrBlank $ws.Cells.Item($row, 9) $antal $true    # I - Antal
    Set-TextOrBlank $ws.Cells.Item($row, 10) $enhet $false  # J - Enhet

    if ($hasL) {
        Set-TextOrBlank $ws.Cells.Item($row, 12) $null $false  # L - Kon (blank, present)
    } else {
        $ws.Cells.Item($row, 12).ClearContents()                # L - Kon (absent)
    }

    $ws.Cells.Item($row, 17).Value2 = $q   # Q - Ost
    $ws.Cells.Item($row, 18).Value2 = $r   # R - Nord
}

# --- capture the current ("before") values for the rows that take part
#     in the rotation, using Value2 which reliably returns scalar data ---

$rows = 2, 3, 5, 6, 7
$before = @{}
foreach ($r in $rows) {
    $before[$r] = @{
        A = $ws.Cells.Item($r, 1).Value2
        B = $ws.Cells.Item($r, 2).Value2
        E = $ws.Cells.Item($r, 5).Value2
        F = $ws.Cells.Item($r, 6).Value2
        G = $ws.Cells.Item($r, 7).Value2
        H = $ws.Cells.Item($r, 8).Value2
        I = $null
        J = $null
        HasL = $false
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
    }
}

# Antal (I) / Enhet (J) are blank for every row in this rotation except row 5.
$before[5].I = "5"
$before[5].J = "fruktkroppar"

# The (empty) placeholder cell in column L is present on rows 3, 6 and 7
# before the edit, and absent on rows 2 and 5.
$before[3].HasL = $true
$before[6].HasL = $true
$before[7].HasL = $true

# --- apply the rotation: row -> gets the data that used to belong to
#     "source" row ---

$sourceOf = @{ 2 = 7; 3 = 5; 5 = 2; 6 = 3; 7 = 6 }

foreach ($r in $rows) {
    $src = $sourceOf[$r]
    $data = $before[$src]

    Set-RowData $r $data.A $data.B $data.E $data.F $data.G $data.H $data.I $data.J $data.HasL $data.Q $data.R
}
